$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting rows 11..89 down to 12..90.
$ws.Rows("11:11").Insert()

# Fill in the new row 11 with the latest weekly price record.
$ws.Range("A11").Value = 11
$ws.Range("B11").Value = "Vega Monumental Concepción"
$ws.Range("C11").Value = "Bíobío"
$ws.Range("D11").Value = 44462
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 100112003
$ws.Range("G11").Value = "Ajo"
$ws.Range("H11").Value = "Chino"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 400
$ws.Range("K11").Value = 15000
$ws.Range("L11").Value = 15500
$ws.Range("M11").Value = 15250
$ws.Range("N11").Value = "$/caja 10 kilos"
$ws.Range("O11").Value = "China"
$ws.Range("P11").Value = 1525
$ws.Range("Q11").Value = 10
$ws.Range("R11").Value = "Hortaliza"
